$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$Bvals = @(0.2330464413454649,0.2063519830744838,0.1899830193775642,0.1833183264427021,0.1822120201958342,0.189893113045116,0.2238379383258007,0.290561330595267,0.3396661954907927,0.3620208745277012,0.3704880823779888,0.3686644363990297,0.3627174395897441,0.3590749773401285,0.3382055637808037,0.3254068183853747,0.3180469086486539,0.3155552554239307,0.3267691051202348,0.3644641663908601,0.3891113930910706,0.3759558053732803,0.3261532206776678,0.2724953414233653)
for ($i = 0; $i -lt $Bvals.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value2 = $Bvals[$i]
}

$Cvals = @(0.0387729706971669,0.03393210848970796,0.03094743649059239,0.02972807408886524,0.02952541432942724,0.03093100420039718,0.03710642065701109,0.04911778533949018,0.05788255978043821,0.06185693635855216,0.06336008013741434,0.06303643473462728,0.0619806384854229,0.06133368881623369,0.05762256544500133,0.05534262047177663,0.05403005695404772,0.05358544043170355,0.05558544887107075,0.06229080233795514,0.06666224522112429,0.06433013092815543,0.05547567173474022,0.04587893463359194)
for ($i = 0; $i -lt $Cvals.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value2 = $Cvals[$i]
}

$Dvals = @(0.02753260113644274,0.02616648665432564,0.02531705957402153,0.02496826440467359,0.02491018812768075,0.02531236628364297,0.02706378552613842,0.03041300423755899,0.03282054332093054,0.03390403039008305,0.0343126102788105,0.03422469193850475,0.03393767900188038,0.03376165133837361,0.03274949672135108,0.03212555192219924,0.03176557446042239,0.03164350406525784,0.03219208605933943,0.03402202838058344,0.03520800262607082,0.0345759502087688,0.03216200990868856,0.02951620108433417)
for ($i = 0; $i -lt $Dvals.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value2 = $Dvals[$i]
}

$Fvals = @(0.7929435651509067,0.7866710901845266,0.7832948023896265,0.7820382717250567,0.7818368298440248,0.7832773733586436,0.7906821462086455,0.8089793609667169,0.8247374550780933,0.8324119115168003,0.8353909834432756,0.834746141264759,0.832655539190327,0.831384486314974,0.8242461070752256,0.819996657446012,0.8176001059694329,0.8167968487747572,0.8204440890426667,0.8332676193596882,0.8420736561632509,0.8373347579000523,0.82024166034968,0.8036238562255562)
for ($i = 0; $i -lt $Fvals.Length; $i++) {
    $ws.Cells.Item($i + 2, 6).Value2 = $Fvals[$i]
}

$Gvals = @(0.6366448308179571,0.6317505978481393,0.6291609484237171,0.628209890116679,0.628058257399033,0.6291477003278203,0.6348709417889609,0.6494018544297262,0.6621122049120345,0.6683403044774252,0.6707631612211316,0.6702384867869853,0.6685383419641511,0.6675053501588764,0.6617141829487707,0.6582759315775633,0.6563403133058472,0.6556921471576374,0.6586375933703437,0.6690359655220561,0.6762074606475466,0.6723454478350845,0.6584739581439436,0.6451148922903371)
for ($i = 0; $i -lt $Gvals.Length; $i++) {
    $ws.Cells.Item($i + 2, 7).Value2 = $Gvals[$i]
}

$Hvals = @(0.7383547964717536,0.7401072477138513,0.7415132626150154,0.7421692041255312,0.7422831343887424,0.7415217729004269,0.7388905462999134,0.7363499345965749,0.7360821516836609,0.7363080312268835,0.7364435956806403,0.7364121739022522,0.7363181813643394,0.7362671243778038,0.7360743890251626,0.7360452265407673,0.7360611778693738,0.7360721983211818,0.736044943826883,0.7363444312330074,0.7368317848238775,0.7365449811616571,0.7360449697374776,0.736756624674868)
for ($i = 0; $i -lt $Hvals.Length; $i++) {
    $ws.Cells.Item($i + 2, 8).Value2 = $Hvals[$i]
}

$Ivals = @(0.5939797809872829,0.5952616649856566,0.596360926302772,0.596887296287143,0.5969794329700271,0.596367707745415,0.5943569273551681,0.5928953752791841,0.5933416391356232,0.5938763653532533,0.5941266750588667,0.5940706376296561,0.5938959993978656,0.5937952596601335,0.5933133787453428,0.5931028061808377,0.5930129060799132,0.5929878251446183,0.5931219904252742,0.5939459960823541,0.5947633226110725,0.5943015506327285,0.5931132201782177,0.5930243520121223)
for ($i = 0; $i -lt $Ivals.Length; $i++) {
    $ws.Cells.Item($i + 2, 9).Value2 = $Ivals[$i]
}

$Kvals = @(0.2460159790660441,0.2159192337892506,0.1974393993125432,0.189908999884409,0.188658610548913,0.1973378400862345,0.2356388935470193,0.3107324393432691,0.3658836301971178,0.3909669490319061,0.4004642990298635,0.3984189327690046,0.391748327275991,0.3876622258820248,0.3642442363190241,0.3498764290503118,0.3416119779160738,0.3388137086698748,0.3514059580768958,0.393707680484539,0.4213473602574425,0.4065963123435949,0.35071447107822,0.2904203632715223)
for ($i = 0; $i -lt $Kvals.Length; $i++) {
    $ws.Cells.Item($i + 2, 11).Value2 = $Kvals[$i]
}

$Mvals = @(0.8860030581400054,0.7889135443472952,0.7298268186455203,0.705874896074647,0.701905157525033,0.7295032907261856,0.8524136059821643,1.097906240168285,1.28141279824986,1.365671473306506,1.397697847486825,1.390794966137179,1.368303866746032,1.354543190385343,1.275922614339521,1.227896433051924,1.200346078842841,1.191030329970118,1.233001291673361,1.374906747997002,1.468349836623418,1.418411058880139,1.230693197072497,1.030972969002676)
for ($i = 0; $i -lt $Mvals.Length; $i++) {
    $ws.Cells.Item($i + 2, 13).Value2 = $Mvals[$i]
}

$Nvals = @(1.527561663344734,1.543899498318957,1.554481742577835,1.558932408031698,1.559679786457096,1.554541206024194,1.533080536983405,1.495373463006807,1.470347812430042,1.459546479170413,1.455540332266516,1.456399383929941,1.459215204418889,1.460950933574424,1.471065457264878,1.477419879673803,1.481129630889122,1.482395102708288,1.476737761397558,1.458385845625997,1.446882032905187,1.45297689320634,1.477045971076102,1.505104673851989)
for ($i = 0; $i -lt $Nvals.Length; $i++) {
    $ws.Cells.Item($i + 2, 14).Value2 = $Nvals[$i]
}
